$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: change formula from D/56+D to D*0.56 ---
# E3 is a standalone formula; E4:E12 is a shared-formula block.
$ws.Range("E3").Formula = "=D3*0.56"
$ws.Range("E4:E12").Formula = "=D4*0.56"

# --- Column G: new column with row averages (B:F) ---
# G3 standalone (matches the pattern already used by column E/F), G4:G12 shared block.
$ws.Range("G3").Formula = "=AVERAGE(B3:F3)"
$ws.Range("G4:G12").Formula = "=AVERAGE(B4:F4)"

# --- Row 14: extend the "PROMEDIO X MES" averages into column G ---
$ws.Range("G14").Formula = "=AVERAGE(G3:G12)"

# --- Row 17: extend the "PRODUCCION TOTAL" array-formula totals into column G ---
$ws.Range("G17").FormulaArray = "=SUM(I14+G3:G12)"

# --- Update the active selection from G3 to H3 ---
$ws.Range("H3").Select()
